# Apply the updated figures to the "Fin Buff Calc" worksheet.
# D3 = Gross Expenditures From 502 Part C
# D5 = Total Labor Cost From 502 Part L
# The dependent formulas in D6, D8, D9, E8, E9 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 117794.69
$ws.Range("D5").Value = 38134.68

$excel.CalculateFullRebuild()
